$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 2-15 (A: empadronador, B: total_registros)
$data = @(
    @{Row=2;  Name="INCIO SANCHEZ PAOLA KATHERINE";      Total=89},
    @{Row=3;  Name="GUEVARA IDROGO DENNIS PERCY";         Total=88},
    @{Row=4;  Name="TANTALEAN BUSTAMANTE ESTALIN YOEL";   Total=86},
    @{Row=5;  Name="HUAYHUA VALDIVIA LUZ EXMILDA";        Total=79},
    @{Row=6;  Name="LINARES PEREZ YANASELY";              Total=77},
    @{Row=7;  Name="MEDINA TAPIA ANA YULI";                Total=77},
    @{Row=8;  Name="PEREZ LINARES TATHIANA";               Total=76},
    @{Row=9;  Name="LOZADA ROJAS LUZ ELENA";               Total=76},
    @{Row=10; Name="CHAVEZ VILLANUEVA SILVIA JANETH";      Total=74},
    @{Row=11; Name="MONDRAGON HERNANDEZ WILMER JUNIOR";    Total=74},
    @{Row=12; Name="DELGADO VASQUEZ FLOR MAGALY";          Total=68},
    @{Row=13; Name="CAMPOS PEREZ YOVERLY";                 Total=65},
    @{Row=14; Name="VASQUEZ SILVA ALOIS ADOLF";            Total=65},
    @{Row=15; Name="SOTO LOZANO LUZDINA";                  Total=62}
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 1).Value = $entry.Name
    $ws.Cells.Item($entry.Row, 2).Value = $entry.Total
}
